$d = $word.ActiveDocument

# Locate the paragraph that holds the "© 2020 ..." footer text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Contact: luizeleno@usp.br*") {
        $target = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the footer paragraph."
}

# The edit removes the footer paragraph together with the two empty
# paragraphs immediately preceding it (a plain spacer paragraph and a
# page-break spacer paragraph), leaving the bibliography entry above
# untouched and the final trailing paragraph intact.
$firstToRemove = $target - 2
$lastToRemove = $target

$startPara = $d.Paragraphs.Item($firstToRemove)
$endPara = $d.Paragraphs.Item($lastToRemove)

$range = $d.Range($startPara.Range.Start, $endPara.Range.End)
$range.Delete()
